# Append a new row (row 6) to the "ランサーズ" sheet and refresh the
# "取得日時" (fetched-at) timestamp on the existing rows to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2025-12-20 12:34:22"

# Update the fetch timestamp on the existing data rows (2-5).
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp

# Append the new record as row 6.
$ws.Range("A6").Value = $newTimestamp
$ws.Range("B6").Value = "Hubspot運用支援(HubSpotのSales Hub・Marketing Hub)"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "300,000 円 ~ 500,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5457876"
$ws.Range("G6").Value = 25

# Turn F6 into a real hyperlink (like F2:F5) and match their "Hyperlink" style.
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5457876")
$ws.Range("F6").Style = $ws.Range("F2").Style
